# "Fixed bug, read formula cell"
#
# Adds a new "FormulaProperty" worksheet (between "IntProperty" and
# "StringProperty") containing a handful of formula cells, so the reader
# bug involving formula cells can be exercised/tested.

$wb = $excel.ActiveWorkbook

# New sheet goes right after "IntProperty" (and therefore right before
# "StringProperty", which keeps its place as the last tab).
$afterSheet = $wb.Worksheets.Item("IntProperty")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = "FormulaProperty"

# Headers
$ws.Range("A1").Value = "Formula"
$ws.Range("F1").Value = "Valor1"
$ws.Range("G1").Value = "Valor2"

# Plain operand values used by the formulas below
$ws.Range("F2").Value = 2
$ws.Range("G3").Value = 4

# Formula cells
$ws.Range("A2").Formula = "=2+3"
$ws.Range("A3").Formula = "=2*4"
$ws.Range("A4").Formula = "=SUM(A2:A3)"
$ws.Range("A5").Formula = "=F2+G3"

# Match the page margins from the saved workbook (metric/A4 defaults).
$ws.PageSetup.LeftMargin = 36.850393728
$ws.PageSetup.RightMargin = 36.850393728
$ws.PageSetup.TopMargin = 56.6929134
$ws.PageSetup.BottomMargin = 56.6929134
$ws.PageSetup.HeaderMargin = 22.67716464
$ws.PageSetup.FooterMargin = 22.67716464

# Selection left where the author's session ended, on the new active sheet.
$ws.Range("A6").Select()
